# Natmi following Dr Hou advice
# Recomputed ligand/receptor/edge stats for the Nid1-Itgav LR pair and
# added the missing "M2" cluster to the sending/target cluster cross-table
# (previously only ECs/FAPs/sCs were present as the 3rd cluster combination).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Nid1"
$ws.Cells.Item(2, 3).Value = "Itgav"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 55.41713066666667
$ws.Cells.Item(2, 8).Value = 166.251392
$ws.Cells.Item(2, 9).Value = 0.08138603925734667
$ws.Cells.Item(2, 10).Value = 0.08138603925734668
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 13.441269
$ws.Cells.Item(2, 14).Value = 40.323807
$ws.Cells.Item(2, 15).Value = 0.08973082133481231
$ws.Cells.Item(2, 16).Value = 0.08973082133481232
$ws.Cells.Item(2, 17).Value = 744.876560498816
$ws.Cells.Item(2, 18).Value = 6703.889044489345
$ws.Cells.Item(2, 19).Value = 0.007302836147748994
$ws.Cells.Item(2, 20).Value = 0.007302836147748996

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Nid1"
$ws.Cells.Item(3, 3).Value = "Itgav"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 55.41713066666667
$ws.Cells.Item(3, 8).Value = 166.251392
$ws.Cells.Item(3, 9).Value = 0.08138603925734667
$ws.Cells.Item(3, 10).Value = 0.08138603925734668
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 54.711535
$ws.Cells.Item(3, 14).Value = 164.134605
$ws.Cells.Item(3, 15).Value = 0.3652416280068742
$ws.Cells.Item(3, 16).Value = 0.3652416280068742
$ws.Cells.Item(3, 17).Value = 3031.956284068907
$ws.Cells.Item(3, 18).Value = 27287.60655662016
$ws.Cells.Item(3, 19).Value = 0.02972556947538467
$ws.Cells.Item(3, 20).Value = 0.02972556947538468

# Row 4: ECs -> M2
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Nid1"
$ws.Cells.Item(4, 3).Value = "Itgav"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 55.41713066666667
$ws.Cells.Item(4, 8).Value = 166.251392
$ws.Cells.Item(4, 9).Value = 0.08138603925734667
$ws.Cells.Item(4, 10).Value = 0.08138603925734668
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 63.67711
$ws.Cells.Item(4, 14).Value = 191.03133
$ws.Cells.Item(4, 15).Value = 0.4250937452800914
$ws.Cells.Item(4, 16).Value = 0.4250937452800915
$ws.Cells.Item(4, 17).Value = 3528.802725345707
$ws.Cells.Item(4, 18).Value = 31759.22452811136
$ws.Cells.Item(4, 19).Value = 0.03459669624141804
$ws.Cells.Item(4, 20).Value = 0.03459669624141806

# Row 5: ECs -> sCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Nid1"
$ws.Cells.Item(5, 3).Value = "Itgav"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 55.41713066666667
$ws.Cells.Item(5, 8).Value = 166.251392
$ws.Cells.Item(5, 9).Value = 0.08138603925734667
$ws.Cells.Item(5, 10).Value = 0.08138603925734668
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 17.96553866666667
$ws.Cells.Item(5, 14).Value = 53.896616
$ws.Cells.Item(5, 15).Value = 0.119933805378222
$ws.Cells.Item(5, 16).Value = 0.119933805378222
$ws.Cells.Item(5, 17).Value = 995.5986037877193
$ws.Cells.Item(5, 18).Value = 8960.387434089473
$ws.Cells.Item(5, 19).Value = 0.00976093739279495
$ws.Cells.Item(5, 20).Value = 0.009760937392794954

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Nid1"
$ws.Cells.Item(6, 3).Value = "Itgav"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 597.374756
$ws.Cells.Item(6, 8).Value = 1792.124268
$ws.Cells.Item(6, 9).Value = 0.8773093221949784
$ws.Cells.Item(6, 10).Value = 0.8773093221949785
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 13.441269
$ws.Cells.Item(6, 14).Value = 40.323807
$ws.Cells.Item(6, 15).Value = 0.08973082133481231
$ws.Cells.Item(6, 16).Value = 0.08973082133481232
$ws.Cells.Item(6, 17).Value = 8029.474789205365
$ws.Cells.Item(6, 18).Value = 72265.27310284828
$ws.Cells.Item(6, 19).Value = 0.07872168604524289
$ws.Cells.Item(6, 20).Value = 0.07872168604524291

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Nid1"
$ws.Cells.Item(7, 3).Value = "Itgav"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 597.374756
$ws.Cells.Item(7, 8).Value = 1792.124268
$ws.Cells.Item(7, 9).Value = 0.8773093221949784
$ws.Cells.Item(7, 10).Value = 0.8773093221949785
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 54.711535
$ws.Cells.Item(7, 14).Value = 164.134605
$ws.Cells.Item(7, 15).Value = 0.3652416280068742
$ws.Cells.Item(7, 16).Value = 0.3652416280068742
$ws.Cells.Item(7, 17).Value = 32683.28987101046
$ws.Cells.Item(7, 18).Value = 294149.6088390941
$ws.Cells.Item(7, 19).Value = 0.3204298851041012
$ws.Cells.Item(7, 20).Value = 0.3204298851041013

# Row 8: FAPs -> M2
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Nid1"
$ws.Cells.Item(8, 3).Value = "Itgav"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 597.374756
$ws.Cells.Item(8, 8).Value = 1792.124268
$ws.Cells.Item(8, 9).Value = 0.8773093221949784
$ws.Cells.Item(8, 10).Value = 0.8773093221949785
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 63.67711
$ws.Cells.Item(8, 14).Value = 191.03133
$ws.Cells.Item(8, 15).Value = 0.4250937452800914
$ws.Cells.Item(8, 16).Value = 0.4250937452800915
$ws.Cells.Item(8, 17).Value = 38039.09804903516
$ws.Cells.Item(8, 18).Value = 342351.8824413164
$ws.Cells.Item(8, 19).Value = 0.3729387055410018
$ws.Cells.Item(8, 20).Value = 0.3729387055410019

# Row 9: FAPs -> sCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Nid1"
$ws.Cells.Item(9, 3).Value = "Itgav"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 597.374756
$ws.Cells.Item(9, 8).Value = 1792.124268
$ws.Cells.Item(9, 9).Value = 0.8773093221949784
$ws.Cells.Item(9, 10).Value = 0.8773093221949785
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 17.96553866666667
$ws.Cells.Item(9, 14).Value = 53.896616
$ws.Cells.Item(9, 15).Value = 0.119933805378222
$ws.Cells.Item(9, 16).Value = 0.119933805378222
$ws.Cells.Item(9, 17).Value = 10732.15927740857
$ws.Cells.Item(9, 18).Value = 96589.4334966771
$ws.Cells.Item(9, 19).Value = 0.1052190455046324
$ws.Cells.Item(9, 20).Value = 0.1052190455046324

# Row 10: M2 -> ECs
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Nid1"
$ws.Cells.Item(10, 3).Value = "Itgav"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.1960536666666667
$ws.Cells.Item(10, 8).Value = 0.5881609999999999
$ws.Cells.Item(10, 9).Value = 0.0002879259755950811
$ws.Cells.Item(10, 10).Value = 0.0002879259755950811
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.441269
$ws.Cells.Item(10, 14).Value = 40.323807
$ws.Cells.Item(10, 15).Value = 0.08973082133481231
$ws.Cells.Item(10, 16).Value = 0.08973082133481232
$ws.Cells.Item(10, 17).Value = 2.635210072103
$ws.Cells.Item(10, 18).Value = 23.716890648927
$ws.Cells.Item(10, 19).Value = 0.00002583583427377375
$ws.Cells.Item(10, 20).Value = 0.00002583583427377376

# Row 11: M2 -> FAPs
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Nid1"
$ws.Cells.Item(11, 3).Value = "Itgav"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1960536666666667
$ws.Cells.Item(11, 8).Value = 0.5881609999999999
$ws.Cells.Item(11, 9).Value = 0.0002879259755950811
$ws.Cells.Item(11, 10).Value = 0.0002879259755950811
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 54.711535
$ws.Cells.Item(11, 14).Value = 164.134605
$ws.Cells.Item(11, 15).Value = 0.3652416280068742
$ws.Cells.Item(11, 16).Value = 0.3652416280068742
$ws.Cells.Item(11, 17).Value = 10.72639704571167
$ws.Cells.Item(11, 18).Value = 96.53757341140499
$ws.Cells.Item(11, 19).Value = 0.0001051625520718149
$ws.Cells.Item(11, 20).Value = 0.000105162552071815

# Row 12: M2 -> M2
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Nid1"
$ws.Cells.Item(12, 3).Value = "Itgav"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.1960536666666667
$ws.Cells.Item(12, 8).Value = 0.5881609999999999
$ws.Cells.Item(12, 9).Value = 0.0002879259755950811
$ws.Cells.Item(12, 10).Value = 0.0002879259755950811
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 63.67711
$ws.Cells.Item(12, 14).Value = 191.03133
$ws.Cells.Item(12, 15).Value = 0.4250937452800914
$ws.Cells.Item(12, 16).Value = 0.4250937452800915
$ws.Cells.Item(12, 17).Value = 12.48413089823667
$ws.Cells.Item(12, 18).Value = 112.35717808413
$ws.Cells.Item(12, 19).Value = 0.0001223955313291372
$ws.Cells.Item(12, 20).Value = 0.0001223955313291373

# Row 13: M2 -> sCs
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Nid1"
$ws.Cells.Item(13, 3).Value = "Itgav"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.1960536666666667
$ws.Cells.Item(13, 8).Value = 0.5881609999999999
$ws.Cells.Item(13, 9).Value = 0.0002879259755950811
$ws.Cells.Item(13, 10).Value = 0.0002879259755950811
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 17.96553866666667
$ws.Cells.Item(13, 14).Value = 53.896616
$ws.Cells.Item(13, 15).Value = 0.119933805378222
$ws.Cells.Item(13, 16).Value = 0.119933805378222
$ws.Cells.Item(13, 17).Value = 3.522209729241778
$ws.Cells.Item(13, 18).Value = 31.699887563176
$ws.Cells.Item(13, 19).Value = 0.00003453205792035515
$ws.Cells.Item(13, 20).Value = 0.00003453205792035517

# Row 14: sCs -> ECs
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Nid1"
$ws.Cells.Item(14, 3).Value = "Itgav"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 27.92897333333333
$ws.Cells.Item(14, 8).Value = 83.78691999999999
$ws.Cells.Item(14, 9).Value = 0.04101671257207978
$ws.Cells.Item(14, 10).Value = 0.04101671257207978
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 13.441269
$ws.Cells.Item(14, 14).Value = 40.323807
$ws.Cells.Item(14, 15).Value = 0.08973082133481231
$ws.Cells.Item(14, 16).Value = 0.08973082133481232
$ws.Cells.Item(14, 17).Value = 375.40084346716
$ws.Cells.Item(14, 18).Value = 3378.60759120444
$ws.Cells.Item(14, 19).Value = 0.00368046330754664
$ws.Cells.Item(14, 20).Value = 0.003680463307546641

# Row 15: sCs -> FAPs
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Nid1"
$ws.Cells.Item(15, 3).Value = "Itgav"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 27.92897333333333
$ws.Cells.Item(15, 8).Value = 83.78691999999999
$ws.Cells.Item(15, 9).Value = 0.04101671257207978
$ws.Cells.Item(15, 10).Value = 0.04101671257207978
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 54.711535
$ws.Cells.Item(15, 14).Value = 164.134605
$ws.Cells.Item(15, 15).Value = 0.3652416280068742
$ws.Cells.Item(15, 16).Value = 0.3652416280068742
$ws.Cells.Item(15, 17).Value = 1528.037002040733
$ws.Cells.Item(15, 18).Value = 13752.3330183666
$ws.Cells.Item(15, 19).Value = 0.01498101087531644
$ws.Cells.Item(15, 20).Value = 0.01498101087531645

# Row 16: sCs -> M2
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Nid1"
$ws.Cells.Item(16, 3).Value = "Itgav"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 27.92897333333333
$ws.Cells.Item(16, 8).Value = 83.78691999999999
$ws.Cells.Item(16, 9).Value = 0.04101671257207978
$ws.Cells.Item(16, 10).Value = 0.04101671257207978
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 63.67711
$ws.Cells.Item(16, 14).Value = 191.03133
$ws.Cells.Item(16, 15).Value = 0.4250937452800914
$ws.Cells.Item(16, 16).Value = 0.4250937452800915
$ws.Cells.Item(16, 17).Value = 1778.436307133733
$ws.Cells.Item(16, 18).Value = 16005.9267642036
$ws.Cells.Item(16, 19).Value = 0.0174359479663424
$ws.Cells.Item(16, 20).Value = 0.01743594796634241

# Row 17: sCs -> sCs
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Nid1"
$ws.Cells.Item(17, 3).Value = "Itgav"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 27.92897333333333
$ws.Cells.Item(17, 8).Value = 83.78691999999999
$ws.Cells.Item(17, 9).Value = 0.04101671257207978
$ws.Cells.Item(17, 10).Value = 0.04101671257207978
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 17.96553866666667
$ws.Cells.Item(17, 14).Value = 53.896616
$ws.Cells.Item(17, 15).Value = 0.119933805378222
$ws.Cells.Item(17, 16).Value = 0.119933805378222
$ws.Cells.Item(17, 17).Value = 501.7590503403022
$ws.Cells.Item(17, 18).Value = 4515.831453062719
$ws.Cells.Item(17, 19).Value = 0.004919290422874287
$ws.Cells.Item(17, 20).Value = 0.004919290422874289
